$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are written as text, not auto-converted to numbers

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.498.16"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.489.79"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.25"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.28"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.488.51"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.944.59"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.354.31"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.11"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.477.10"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.89"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.87"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.89"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.21"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.615.42"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.58"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0865"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.54"
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "435.05"
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "154.98"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.07"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.57"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.40"
$ws.Range("E44").Value = "  +50.27%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "138.22"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.503"
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0723"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.572"
$ws.Range("E51").Value = "  -0.82%  "
